$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 97.25956416130066
$ws.Range("C2").Value = 0.7527976803821724
$ws.Range("D2").Value = 0.9662322044372559
$ws.Range("E2").Value = 0.02312501834607601
$ws.Range("J2").Value = 0.7017543859649122
$ws.Range("K2").Value = 0.6081871345029239
$ws.Range("L2").Value = 0.695906432748538
$ws.Range("M2").Value = 0.6900584795321637
$ws.Range("N2").Value = 0.6549707602339181
$ws.Range("O2").Value = 0.6701754385964913
$ws.Range("P2").Value = 0.03500965975928835
$ws.Range("Q2").Value = 13

# Row 3
$ws.Range("B3").Value = 193.2165286064148
$ws.Range("C3").Value = 1.584161975291451
$ws.Range("D3").Value = 1.010275411605835
$ws.Range("E3").Value = 0.1234121021173341
$ws.Range("J3").Value = 0.7251461988304093
$ws.Range("K3").Value = 0.6198830409356725
$ws.Range("L3").Value = 0.6900584795321637
$ws.Range("M3").Value = 0.6549707602339181
$ws.Range("N3").Value = 0.672514619883041
$ws.Range("O3").Value = 0.672514619883041
$ws.Range("P3").Value = 0.03508771929824561
$ws.Range("Q3").Value = 10

# Row 4
$ws.Range("B4").Value = 371.7606033802032
$ws.Range("C4").Value = 3.168798320844412
$ws.Range("D4").Value = 1.021530246734619
$ws.Range("E4").Value = 0.1290915228128792
$ws.Range("J4").Value = 0.7192982456140351
$ws.Range("K4").Value = 0.5964912280701754
$ws.Range("L4").Value = 0.7017543859649122
$ws.Range("M4").Value = 0.6549707602339181
$ws.Range("N4").Value = 0.7017543859649122
$ws.Range("O4").Value = 0.6748538011695906
$ws.Range("P4").Value = 0.04462873456921385
$ws.Range("Q4").Value = 4

# Row 5
$ws.Range("B5").Value = 95.95504670143127
$ws.Range("C5").Value = 1.253701205062879
$ws.Range("D5").Value = 0.9608968734741211
$ws.Range("E5").Value = 0.0140618445553067
$ws.Range("J5").Value = 0.7017543859649122
$ws.Range("K5").Value = 0.6081871345029239
$ws.Range("L5").Value = 0.695906432748538
$ws.Range("M5").Value = 0.6900584795321637
$ws.Range("N5").Value = 0.6549707602339181
$ws.Range("O5").Value = 0.6701754385964913
$ws.Range("P5").Value = 0.03500965975928835
$ws.Range("Q5").Value = 13

# Row 6
$ws.Range("B6").Value = 187.6987954139709
$ws.Range("C6").Value = 1.01061377056972
$ws.Range("D6").Value = 0.9766998291015625
$ws.Range("E6").Value = 0.03037780577802829
$ws.Range("J6").Value = 0.7251461988304093
$ws.Range("K6").Value = 0.6198830409356725
$ws.Range("L6").Value = 0.6900584795321637
$ws.Range("M6").Value = 0.6549707602339181
$ws.Range("N6").Value = 0.672514619883041
$ws.Range("O6").Value = 0.672514619883041
$ws.Range("P6").Value = 0.03508771929824561
$ws.Range("Q6").Value = 10

# Row 7
$ws.Range("B7").Value = 368.9238994598389
$ws.Range("C7").Value = 1.319249664022278
$ws.Range("D7").Value = 0.9437858581542968
$ws.Range("E7").Value = 0.003731690985331589
$ws.Range("J7").Value = 0.7192982456140351
$ws.Range("K7").Value = 0.5964912280701754
$ws.Range("L7").Value = 0.7017543859649122
$ws.Range("M7").Value = 0.6549707602339181
$ws.Range("N7").Value = 0.7017543859649122
$ws.Range("O7").Value = 0.6748538011695906
$ws.Range("P7").Value = 0.04462873456921385
$ws.Range("Q7").Value = 4

# Row 8
$ws.Range("B8").Value = 96.0635157585144
$ws.Range("C8").Value = 1.007650687502769
$ws.Range("D8").Value = 1.033104848861694
$ws.Range("E8").Value = 0.1274844277002802
$ws.Range("J8").Value = 0.7017543859649122
$ws.Range("K8").Value = 0.6081871345029239
$ws.Range("L8").Value = 0.695906432748538
$ws.Range("M8").Value = 0.6900584795321637
$ws.Range("N8").Value = 0.6549707602339181
$ws.Range("O8").Value = 0.6701754385964913
$ws.Range("P8").Value = 0.03500965975928835
$ws.Range("Q8").Value = 13

# Row 9
$ws.Range("B9").Value = 184.7123771190643
$ws.Range("C9").Value = 2.22242808129428
$ws.Range("D9").Value = 1.031090116500855
$ws.Range("E9").Value = 0.1271741464987291
$ws.Range("J9").Value = 0.7251461988304093
$ws.Range("K9").Value = 0.6198830409356725
$ws.Range("L9").Value = 0.6900584795321637
$ws.Range("M9").Value = 0.6549707602339181
$ws.Range("N9").Value = 0.672514619883041
$ws.Range("O9").Value = 0.672514619883041
$ws.Range("P9").Value = 0.03508771929824561
$ws.Range("Q9").Value = 10

# Row 10
$ws.Range("B10").Value = 362.1925025939942
$ws.Range("C10").Value = 2.034981903479518
$ws.Range("D10").Value = 0.986182689666748
$ws.Range("E10").Value = 0.03550121344016233
$ws.Range("J10").Value = 0.7192982456140351
$ws.Range("K10").Value = 0.5964912280701754
$ws.Range("L10").Value = 0.7017543859649122
$ws.Range("M10").Value = 0.6549707602339181
$ws.Range("N10").Value = 0.7017543859649122
$ws.Range("O10").Value = 0.6748538011695906
$ws.Range("P10").Value = 0.04462873456921385
$ws.Range("Q10").Value = 4

# Row 11
$ws.Range("B11").Value = 144.7844655036926
$ws.Range("C11").Value = 1.382558774848798
$ws.Range("D11").Value = 0.9461830139160157
$ws.Range("E11").Value = 0.02778496890510089
$ws.Range("J11").Value = 0.7251461988304093
$ws.Range("K11").Value = 0.6081871345029239
$ws.Range("L11").Value = 0.6608187134502924
$ws.Range("M11").Value = 0.7017543859649122
$ws.Range("N11").Value = 0.6783625730994152
$ws.Range("O11").Value = 0.6748538011695906
$ws.Range("P11").Value = 0.03976608187134503
$ws.Range("Q11").Value = 4

# Row 12
$ws.Range("B12").Value = 280.0421406269073
$ws.Range("C12").Value = 0.9048531301649925
$ws.Range("D12").Value = 1.173026323318481
$ws.Range("E12").Value = 0.1239467541938699
$ws.Range("J12").Value = 0.7134502923976608
$ws.Range("K12").Value = 0.5847953216374269
$ws.Range("L12").Value = 0.6842105263157895
$ws.Range("M12").Value = 0.7076023391812866
$ws.Range("N12").Value = 0.6549707602339181
$ws.Range("O12").Value = 0.6690058479532164
$ws.Range("P12").Value = 0.04687126294623131

# Row 13
$ws.Range("B13").Value = 540.5873781681061
$ws.Range("C13").Value = 3.080923152560644
$ws.Range("D13").Value = 0.9898348808288574
$ws.Range("E13").Value = 0.08283058639919517
$ws.Range("J13").Value = 0.7309941520467836
$ws.Range("K13").Value = 0.5906432748538012
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.7251461988304093
$ws.Range("N13").Value = 0.6666666666666666
$ws.Range("O13").Value = 0.6760233918128654
$ws.Range("P13").Value = 0.05079310117101394
$ws.Range("Q13").Value = 1

# Row 14
$ws.Range("B14").Value = 143.9262855529785
$ws.Range("C14").Value = 1.491898179360678
$ws.Range("D14").Value = 1.075334930419922
$ws.Range("E14").Value = 0.1803250783778788
$ws.Range("J14").Value = 0.7251461988304093
$ws.Range("K14").Value = 0.6081871345029239
$ws.Range("L14").Value = 0.6608187134502924
$ws.Range("M14").Value = 0.7017543859649122
$ws.Range("N14").Value = 0.6783625730994152
$ws.Range("O14").Value = 0.6748538011695906
$ws.Range("P14").Value = 0.03976608187134503
$ws.Range("Q14").Value = 4

# Row 15
$ws.Range("B15").Value = 280.9537870407104
$ws.Range("C15").Value = 0.7425645618745271
$ws.Range("D15").Value = 0.9150830268859863
$ws.Range("E15").Value = 0.009408836164922503
$ws.Range("J15").Value = 0.7134502923976608
$ws.Range("K15").Value = 0.5847953216374269
$ws.Range("L15").Value = 0.6842105263157895
$ws.Range("M15").Value = 0.7076023391812866
$ws.Range("N15").Value = 0.6549707602339181
$ws.Range("O15").Value = 0.6690058479532164
$ws.Range("P15").Value = 0.04687126294623131

# Row 16
$ws.Range("B16").Value = 542.6661983966827
$ws.Range("C16").Value = 2.605783449308912
$ws.Range("D16").Value = 0.9872122764587402
$ws.Range("E16").Value = 0.0711649755272738
$ws.Range("J16").Value = 0.7309941520467836
$ws.Range("K16").Value = 0.5906432748538012
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.7251461988304093
$ws.Range("N16").Value = 0.6666666666666666
$ws.Range("O16").Value = 0.6760233918128654
$ws.Range("P16").Value = 0.05079310117101394
$ws.Range("Q16").Value = 1

# Row 17
$ws.Range("B17").Value = 143.088559961319
$ws.Range("C17").Value = 0.6988840550255948
$ws.Range("D17").Value = 1.01957745552063
$ws.Range("E17").Value = 0.1249629178655173
$ws.Range("J17").Value = 0.7251461988304093
$ws.Range("K17").Value = 0.6081871345029239
$ws.Range("L17").Value = 0.6608187134502924
$ws.Range("M17").Value = 0.7017543859649122
$ws.Range("N17").Value = 0.6783625730994152
$ws.Range("O17").Value = 0.6748538011695906
$ws.Range("P17").Value = 0.03976608187134503
$ws.Range("Q17").Value = 4

# Row 18
$ws.Range("B18").Value = 281.4950689315796
$ws.Range("C18").Value = 1.718722566403527
$ws.Range("D18").Value = 0.9799116134643555
$ws.Range("E18").Value = 0.09180529056139601
$ws.Range("J18").Value = 0.7134502923976608
$ws.Range("K18").Value = 0.5847953216374269
$ws.Range("L18").Value = 0.6842105263157895
$ws.Range("M18").Value = 0.7076023391812866
$ws.Range("N18").Value = 0.6549707602339181
$ws.Range("O18").Value = 0.6690058479532164
$ws.Range("P18").Value = 0.04687126294623131

# Row 19
$ws.Range("B19").Value = 455.4692974090576
$ws.Range("C19").Value = 19.18992232435417
$ws.Range("D19").Value = 0.6252786159515381
$ws.Range("E19").Value = 0.1170188300016989
$ws.Range("J19").Value = 0.7309941520467836
$ws.Range("K19").Value = 0.5906432748538012
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.7251461988304093
$ws.Range("N19").Value = 0.6666666666666666
$ws.Range("O19").Value = 0.6760233918128654
$ws.Range("P19").Value = 0.05079310117101394
$ws.Range("Q19").Value = 1
